$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values below are plain strings in the source data (inline strings in the
# canonical OOXML). A handful of the new Price values parse as ordinary decimals
# (e.g. "1.00", "0.101"), so Excel would silently coerce them to numbers -- and drop
# the literal formatting (trailing zeros, leading zero, etc.) -- if assigned directly.
# For those specific cells we briefly force a Text number format so the assignment
# is kept verbatim, then restore the default "Normal" style afterwards.

$ws.Range('D2').Value = '63.180.01'
$ws.Range('E2').Value = '  -7.09%  '
$ws.Range('D3').Value = '3.228.90'
$ws.Range('E3').Value = '  -9.84%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '173.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -16.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '507.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -11.04%  '
$ws.Range('E7').Value = '  -4.77%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '3.219.27'
$ws.Range('E9').Value = '  -9.86%  '
$ws.Range('E10').Value = '  -10.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.23%  '
$ws.Range('E13').Value = '  -10.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.91%  '
$ws.Range('D15').Value = '3.751.88'
$ws.Range('E15').Value = '  -9.56%  '
$ws.Range('E16').Value = '  -6.93%  '
$ws.Range('D17').Value = '3.232.56'
$ws.Range('E17').Value = '  -9.68%  '
$ws.Range('D18').Value = '62.997.21'
$ws.Range('E18').Value = '  -7.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.94%  '
$ws.Range('E20').Value = '  -13.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.926'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -12.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '363.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '78.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.29%  '
$ws.Range('E24').Value = '  -14.39%  '
$ws.Range('E25').Value = '  -14.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.15%  '
$ws.Range('E28').Value = '  -10.15%  '
$ws.Range('E29').Value = '  -11.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '634.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.24%  '
$ws.Range('E32').Value = '  -11.69%  '
$ws.Range('E33').Value = '  -15.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.53%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.101'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -10.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '34.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -14.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.369'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.79%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  -9.22%  '
$ws.Range('D42').Value = '2.831.16'
$ws.Range('E42').Value = '  -10.69%  '
$ws.Range('E43').Value = '  -14.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -18.82%  '
$ws.Range('E45').Value = '  -7.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -13.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.73'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0373'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.34%  '
$ws.Range('E50').Value = '  -6.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '130.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.58%  '
